$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G12").Value = 1.72
$ws.Range("H12").Value = 3.6
$ws.Range("J12").Value = 2.25
$ws.Range("K12").Value = 2.18
$ws.Range("L12").Value = 4.6
$ws.Range("U12").Value = 1.62
$ws.Range("V12").Value = 2.02
$ws.Range("X12").Value = 9
$ws.Range("Z12").Value = 14.5
$ws.Range("AD12").Value = 7.1
$ws.Range("AE12").Value = 13.5
$ws.Range("AI12").Value = 14
$ws.Range("AM12").Value = 350
$ws.Range("AN12").Value = 3.65
$ws.Range("AO12").Value = 8.25
$ws.Range("AQ12").Value = 27
$ws.Range("AT12").Value = 2.87
$ws.Range("AU12").Value = 6.9
$ws.Range("AV12").Value = 55
$ws.Range("AW12").Value = 6.2
$ws.Range("AX12").Value = 24
$ws.Range("N16").Value = 9
$ws.Range("Q16").Value = 2.2
$ws.Range("R16").Value = 1.65
$ws.Range("M17").Value = 1.06
$ws.Range("N17").Value = 10
$ws.Range("O17").Value = 1.33
$ws.Range("P17").Value = 3.25
$ws.Range("I18").Value = 2.55
$ws.Range("J18").Value = 3.2
$ws.Range("K18").Value = 2.1
$ws.Range("N18").Value = 6.9
$ws.Range("O18").Value = 1.34
$ws.Range("P18").Value = 3
$ws.Range("Q18").Value = 2.05
$ws.Range("R18").Value = 1.72
$ws.Range("U18").Value = 1.82
$ws.Range("V18").Value = 1.9
$ws.Range("W18").Value = 7.9
$ws.Range("X18").Value = 12.5
$ws.Range("AA18").Value = 22
$ws.Range("AC18").Value = 6.9
$ws.Range("AE18").Value = 15
$ws.Range("AF18").Value = 75
$ws.Range("AG18").Value = 7.9
$ws.Range("AH18").Value = 12
$ws.Range("AJ18").Value = 27
$ws.Range("AL18").Value = 32
$ws.Range("AM18").Value = 600
$ws.Range("AO18").Value = 14
$ws.Range("AW18").Value = 4.45
$ws.Range("G28").Value = 2.57
$ws.Range("H28").Value = 3.2
$ws.Range("I28").Value = 2.52
$ws.Range("J28").Value = 3.25
$ws.Range("L28").Value = 3.15
$ws.Range("N28").Value = 7.2
$ws.Range("O28").Value = 1.3
$ws.Range("P28").Value = 3.2
$ws.Range("Q28").Value = 1.91
$ws.Range("R28").Value = 1.83
$ws.Range("T28").Value = 2.65
$ws.Range("X28").Value = 13
$ws.Range("Y28").Value = 9.75
$ws.Range("Z28").Value = 29
$ws.Range("AA28").Value = 22
$ws.Range("AB28").Value = 30
$ws.Range("AC28").Value = 7.2
$ws.Range("AD28").Value = 6.3
$ws.Range("AG28").Value = 8.75
$ws.Range("AH28").Value = 13
$ws.Range("AI28").Value = 9.5
$ws.Range("AJ28").Value = 28
$ws.Range("AK28").Value = 21
$ws.Range("AL28").Value = 29
$ws.Range("AN28").Value = 4.55
$ws.Range("AO28").Value = 14.5
$ws.Range("AP28").Value = 22
$ws.Range("AQ28").Value = 65
$ws.Range("AR28").Value = 100
$ws.Range("AT28").Value = 2.65
$ws.Range("AU28").Value = 7
$ws.Range("AV28").Value = 65
$ws.Range("AW28").Value = 4.5
$ws.Range("AX28").Value = 14
$ws.Range("AY28").Value = 22
$ws.Range("AZ28").Value = 60
$ws.Range("BA28").Value = 90
$ws.Range("BB28").Value = 300
$ws.Range("G29").Value = 2.52
$ws.Range("H29").Value = 3.5
$ws.Range("I29").Value = 2.4
$ws.Range("K29").Value = 2.2
$ws.Range("L29").Value = 3
$ws.Range("N29").Value = 7.9
$ws.Range("O29").Value = 1.26
$ws.Range("R29").Value = 1.95
$ws.Range("S29").Value = 1.37
$ws.Range("T29").Value = 2.87
$ws.Range("AC29").Value = 7.9
$ws.Range("AD29").Value = 6.9
$ws.Range("AH29").Value = 12.5
$ws.Range("AJ29").Value = 25
$ws.Range("AK29").Value = 19
$ws.Range("AT29").Value = 2.87
$ws.Range("AU29").Value = 7.1
$ws.Range("AW29").Value = 4.45
$ws.Range("AX29").Value = 12.5
$ws.Range("AZ29").Value = 50
$ws.Range("BA29").Value = 80
$ws.Range("L30").Value = 4.8
$ws.Range("P30").Value = 5.1
$ws.Range("Q30").Value = 1.42
$ws.Range("R30").Value = 2.67
$ws.Range("AM30").Value = 250
$ws.Range("AQ30").Value = 17.5
$ws.Range("AU30").Value = 7
$ws.Range("AW30").Value = 7.3
$ws.Range("G31").Value = 1.72
$ws.Range("I31").Value = 4.1
$ws.Range("L31").Value = 4.5
$ws.Range("Q31").Value = 1.83
$ws.Range("T31").Value = 2.8
$ws.Range("W31").Value = 7.3
$ws.Range("X31").Value = 8.25
$ws.Range("AI31").Value = 14
$ws.Range("AP31").Value = 18
$ws.Range("AQ31").Value = 30
$ws.Range("AT31").Value = 2.8
$ws.Range("AU31").Value = 7.5
$ws.Range("AV31").Value = 70
$ws.Range("AW31").Value = 6
$ws.Range("G32").Value = 2.57
$ws.Range("L32").Value = 3.15
$ws.Range("R32").Value = 1.85
$ws.Range("AW32").Value = 4.55
$ws.Range("H33").Value = 3.65
$ws.Range("I33").Value = 1.95
$ws.Range("K33").Value = 2.22
$ws.Range("P33").Value = 3.7
$ws.Range("S33").Value = 1.36
$ws.Range("T33").Value = 2.92
$ws.Range("W33").Value = 11.5
$ws.Range("Y33").Value = 11.5
$ws.Range("AL33").Value = 23
$ws.Range("AN33").Value = 5.2
$ws.Range("AO33").Value = 17.5
$ws.Range("AP33").Value = 24
$ws.Range("AT33").Value = 2.92
$ws.Range("AY33").Value = 18
$ws.Range("G38").Value = 10.5
$ws.Range("H38").Value = 5.1
$ws.Range("O38").Value = 1.22
$ws.Range("S38").Value = 1.33
$ws.Range("T38").Value = 3.05
$ws.Range("U38").Value = 2.32
$ws.Range("W38").Value = 25
$ws.Range("X38").Value = 80
$ws.Range("Y38").Value = 35
$ws.Range("Z38").Value = 400
$ws.Range("AD38").Value = 11
$ws.Range("AE38").Value = 30
$ws.Range("AJ38").Value = 6.8
$ws.Range("AK38").Value = 11.75
$ws.Range("AN38").Value = 10.75
$ws.Range("AO38").Value = 70
$ws.Range("AP38").Value = 65
$ws.Range("AT38").Value = 3.05
$ws.Range("H39").Value = 4.45
$ws.Range("I39").Value = 6.4
$ws.Range("K39").Value = 2.47
$ws.Range("L39").Value = 5.9
$ws.Range("U39").Value = 1.78
$ws.Range("V39").Value = 1.93
$ws.Range("W39").Value = 8.25
$ws.Range("X39").Value = 7.4
$ws.Range("Z39").Value = 9.5
$ws.Range("AB39").Value = 23
$ws.Range("AG39").Value = 21
$ws.Range("AQ39").Value = 17
$ws.Range("AV39").Value = 60
$ws.Range("AZ39").Value = 200
$ws.Range("BA39").Value = 200
$ws.Range("J40").Value = 2.4
$ws.Range("U40").Value = 1.52
$ws.Range("V40").Value = 2.37
$ws.Range("X40").Value = 11.5
$ws.Range("AB40").Value = 19
$ws.Range("AE40").Value = 12
$ws.Range("AM40").Value = 250
$ws.Range("AO40").Value = 9.25
$ws.Range("AR40").Value = 45
$ws.Range("G41").Value = 3.7
$ws.Range("H41").Value = 3.35
$ws.Range("J41").Value = 4.15
$ws.Range("K41").Value = 2.1
$ws.Range("N41").Value = 7.4
$ws.Range("Q41").Value = 1.87
$ws.Range("U41").Value = 1.72
$ws.Range("V41").Value = 2
$ws.Range("W41").Value = 11.75
$ws.Range("X41").Value = 22
$ws.Range("AC41").Value = 7.4
$ws.Range("AD41").Value = 6.6
$ws.Range("AG41").Value = 7.3
$ws.Range("AH41").Value = 9.25
$ws.Range("AL41").Value = 25
$ws.Range("AP41").Value = 26
$ws.Range("AU41").Value = 7.1
$ws.Range("AX41").Value = 10
$ws.Range("I42").Value = 3.35
$ws.Range("R42").Value = 1.82
$ws.Range("V42").Value = 1.95
$ws.Range("AD42").Value = 6.5
$ws.Range("AH42").Value = 18.5
$ws.Range("AJ42").Value = 45
$ws.Range("AU42").Value = 7.2
$ws.Range("AX42").Value = 18.5
$ws.Range("AZ42").Value = 90
$ws.Range("BB42").Value = 350
